# Revert config file handling
# Appends a new row (row 45) of decoded/raw config-packet data to each of
# the four sheets (MID_LFT_#1, MID_LFT_#2, MID_PLT_#1, MID_PLT_#2), mirroring
# the structure of the existing rows (1..44) in every sheet.

$wb = $excel.ActiveWorkbook

# Per-sheet new-row values, in column order: A,B,C,D,E,F,G,H,I
$rowsBySheet = @{
    "MID_LFT_#1" = @(
        45831.46317129629,
        "0x01,0x90",
        "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,",
        "0x01,0x70",
        "0x07",
        400,
        [double]"5.68631262647113e+23",
        368,
        7
    )
    "MID_LFT_#2" = @(
        45831.46317129629,
        "0x01,0x7c",
        "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,",
        "0x01,0x64",
        "0x19",
        380,
        [double]"5.68432987514711e+23",
        356,
        25
    )
    "MID_PLT_#1" = @(
        45831.46317129629,
        "0x00,0x6e",
        "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,",
        "0x00,0x6A",
        "0x15",
        110,
        [double]"5.68631262647113e+23",
        106,
        15
    )
    "MID_PLT_#2" = @(
        45831.46317129629,
        "0x00,0x82",
        "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,",
        "0x00,0x7F",
        "0x9",
        130,
        [double]"5.68631262647113e+23",
        127,
        9
    )
}

$newRow = 45
$sourceRow = 44

foreach ($ws in $wb.Worksheets) {
    $values = $rowsBySheet[$ws.Name]
    if ($null -eq $values) {
        continue
    }

    for ($col = 1; $col -le 9; $col++) {
        $dstCell = $ws.Cells.Item($newRow, $col)
        $dstCell.Value = $values[$col - 1]
    }

    # Column A carries the date/time number format - match the row above it.
    # (Columns B..I use the default "General" style, same as every other
    # data row, so they are left untouched.)
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($sourceRow, 1).NumberFormat
}
